# Amendment name, label and description change, DDF-RA #625
#
# studyAmendments sheet: insert three new leading columns (name, label,
# description) ahead of the existing number/summary/... columns, and
# populate them with the amendment identifiers/labels. Also makes the
# studyAmendments sheet the active/selected sheet (previously it was
# studyDesignPopulations).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyAmendments")

# --- Insert 3 new columns at the front (A:C), shifting old A:G to D:J ---
$ws.Columns("A:C").Insert()

# --- Header row (row 1): reuse the header look (fill/font) from the
#     existing header cell, then force left/top alignment to match the
#     "name"/"label"/"description" header style used elsewhere. ---
$ws.Range("D1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Range("A1:C1").HorizontalAlignment = -4131

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "label"
$ws.Range("C1").Value = "description"

# --- Data rows: plain left-aligned cells ---
$ws.Range("A2:C5").HorizontalAlignment = -4131

$ws.Range("A2").Value = "AMEND_1"
$ws.Range("A3").Value = "AMEND_2"
$ws.Range("A4").Value = "AMEND_3"
$ws.Range("A5").Value = "AMEND_4"

$ws.Range("B2").Value = "Amendment 1"
$ws.Range("B3").Value = "Amendment 2"
$ws.Range("B4").Value = "Amendment 3"
$ws.Range("B5").Value = "Amendment 4"

$ws.Range("C2").Value = "Amendment 1"
$ws.Range("C3").Value = "Amendment 2"
$ws.Range("C4").Value = "Amendment 3"
$ws.Range("C5").Value = "Amendment 4"

# --- Column widths for the two newly inserted, width-bearing columns ---
$ws.Columns("B").ColumnWidth = 16.1666666667
$ws.Columns("C").ColumnWidth = 14.1666666667

# --- Make studyAmendments the active sheet / tab, with D13 selected
#     (previously studyDesignPopulations was the selected tab). ---
$ws.Activate()
$ws.Range("D13").Select()
